$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "2022" column (S) -----------------------------------
# Year header S4: clone formatting from the adjacent year header R4,
# then set its value.
$ws.Range("R4").Copy() | Out-Null
$ws.Range("S4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("S4").Value = 2022

# Data rows 5-13: clone formatting from the default data style (column S
# picks up the plain default style automatically, matching column style 4),
# so a plain value assignment is enough.
$ws.Range("S5").Value = 115.8
$ws.Range("S6").Value = 115.2
$ws.Range("S7").Value = 115.4
$ws.Range("S8").Value = 111.8
$ws.Range("S9").Value = 116.8
$ws.Range("S10").Value = 108.2
$ws.Range("S11").Value = 111
$ws.Range("S12").Value = 115.8
$ws.Range("S13").Value = 117.9

# Row 14 (totals row) needs the heavier bottom-border style used by the
# rest of that row (matching R14), so clone formatting first.
$ws.Range("R14").Copy() | Out-Null
$ws.Range("S14").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("S14").Value = 112.4

$excel.CutCopyMode = $false

# --- Refresh the 2020/2021 figures with the restated series -----------
$ws.Range("Q5").Value = 117.60684979252385
$ws.Range("R5").Value = 113.34848864817617

$ws.Range("Q6").Value = 114.77319768114526
$ws.Range("R6").Value = 115.06069350712495

$ws.Range("Q7").Value = 116.40044011407315
$ws.Range("R7").Value = 114.29658549692938

$ws.Range("Q8").Value = 117.53828537152096
$ws.Range("R8").Value = 113.75761785228545

$ws.Range("Q9").Value = 117.42206669681742
$ws.Range("R9").Value = 113.98264089946031

$ws.Range("Q10").Value = 113.98326995089161
$ws.Range("R10").Value = 113.92720567782911

$ws.Range("Q11").Value = 123.488978736909
$ws.Range("R11").Value = 114.17226706705155

$ws.Range("Q12").Value = 118.12340252754679
$ws.Range("R12").Value = 114.45153946490467

$ws.Range("Q13").Value = 118.87059844457349
$ws.Range("R13").Value = 112.69493421065988

$ws.Range("Q14").Value = 114.06377070452145
$ws.Range("R14").Value = 113.95067699644588

# --- Move the active-cell selection from T6 to T4 ----------------------
$ws.Range("T4").Select() | Out-Null

Write-Host "2.c.1.1b updated: added 2022 column, refreshed 2020/2021 figures"
